# Update crypto price/volume data (Price and Volume(1h) columns) on Sheet1
# to reflect the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.143.66'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.420.59'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'553.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = "'136.96"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').Value = "'24.89"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').Value = '2.850.44'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '60.045.51'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').Value = '2.408.25'
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = "'11.29"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('E19').Value = '  +1.54%  '
$ws.Range('D20').Value = "'328.24"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = "'0.176"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('D25').Value = "'8.71"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = "'1.38"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.19%  '
$ws.Range('E28').Value = '  -2.67%  '
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('D30').Value = "'170.09"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('E32').Value = '  +2.06%  '
$ws.Range('D33').Value = "'0.404"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').Value = "'328.14"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.20%  '
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('D41').Value = "'38.68"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('D42').Value = "'145.28"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.64%  '
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').Value = "'20.01"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.11%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = "'11.04"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('E51').Value = '  -0.85%  '
